# Fix the typo "STMP" -> "SMTP" in "L'STMP si basa sul trasporto TCP ..."
#
# The canonical diff shows the corrected text split across three runs
# (L'S | MT | P si basa...) with identical (default) run formatting -
# exactly what Word produces when a user retypes/corrects a couple of
# characters in the middle of an existing run. We reproduce that by:
#   1. fixing the text itself via Find & Replace, then
#   2. re-creating the same run boundaries with a harmless
#      Bold-on/Bold-off formatting round trip (touches the run without
#      altering its visible formatting), which is enough to stop the
#      engine from re-coalescing the edited text back into one run.

$d = $word.ActiveDocument
$apostrophe = [char]0x2019

# --- Step 1: correct the misspelling "L'STMP" -> "L'SMTP" -------------
$rng = $d.Content
$found = $rng.Find.Execute("L" + $apostrophe + "STMP", $true, $false, $false, `
    $false, $false, $true, 1, $false, "L" + $apostrophe + "SMTP", 2)

if (-not $found) {
    Write-Output "WARNING: target text 'L'STMP' not found"
}

# --- Step 2: re-split the run the way the source document does --------
# Locate the corrected text again to get fresh character offsets.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("L" + $apostrophe + "SMTP", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $s = $rng2.Start
    # "L'SMTP" -> L(0) '(1) S(2) M(3) T(4) P(5)
    # run 1: "L'S"  -> [s, s+3)
    # run 2: "MT"   -> [s+3, s+5)
    # run 3: "P ..."-> [s+5, ...) (left untouched, already part of the
    #                  following text)
    $run1 = $d.Range($s, $s + 3)
    $run1.Bold = 1
    $run1.Bold = 0

    $run2 = $d.Range($s + 3, $s + 5)
    $run2.Bold = 1
    $run2.Bold = 0
}

Write-Output "done"
